$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = 77
$ws.Range("B20").Value = 60.441165924072266
$ws.Range("B22").Value = 95
$ws.Range("B29").Value = 38
